$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Values that are purely numeric-looking
# are apostrophe-prefixed so Excel stores them as Text (matching the source
# workbook, where every cell -- including price figures -- is inline text).
$updates = [ordered]@{
    "D2" = "258.07"
    "D3" = "21.42"
    "D4" = "6.126"
    "D5" = "0.06131"
    "D6" = "3.572"
    "D7" = "6.538"
    "D8" = "1.371"
    "D9" = "0.8239"
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D10" = "0.01323"
    "E10" = "9OneONE"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "0.1609"
    "E11" = "10WazirXWRX"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "0.08139"
    "E12" = "11MandalaExchangeTokenMDX"
    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D13" = "0.03530"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D14" = "0.03175"
    "E14" = "13BitrueCoinBTR"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D15" = "0.09221"
    "E15" = "14BitMartTokenBMX"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D16" = "3.782"
    "E16" = "15MCDexMCB"
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D17" = "0.001641"
    "E17" = "16BitForexTokenBF"
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D18" = "0.04653"
    "E18" = "17CoinExTokenCET"
    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D19" = "0.006422"
    "E19" = "18TigerCashTCH"
    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "D20" = "0.006133"
    "E20" = "19HotbitTokenHTB"
    "B21" = "BitKan"
    "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "D21" = "0.001070"
    "E21" = "20BitKanKAN"
    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "D22" = "0.0001505"
    "E22" = "21NitroExNTX"
    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D23" = "3.725"
    "E23" = "22LEOLEO"
    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D24" = "2.271"
    "E24" = "23BTSETokenBTSE"
    "D25" = "0.3317"
    "D28" = "0.0002722"
    "D40" = "0.04666"
    "D41" = "0.007013"
    "D42" = "0.003732"
    "D43" = "0.1115"
    "D45" = "0.00006062"
    "D47" = "0.00000000752"
    "D48" = "0.9832"
    "D49" = "0.001132"
    "D50" = "0.00001906"
    "D51" = "0.01244"
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $looksNumeric = $value -match '^[0-9]*\.?[0-9]+$'
    if ($looksNumeric) {
        # Apostrophe-prefix forces Excel to keep/store this as Text instead
        # of silently coercing the numeric-looking string to a Number.
        $ws.Range($addr).Value = "'" + $value
    } else {
        $ws.Range($addr).Value = $value
    }
}

Write-Host "Applied $($updates.Count) cell updates"
